# Auto-generated Excel COM-interop script to apply scheduled-runner price/profit updates
# to the Spriggan_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 144610.5
$ws.Range("J17").Value = 144610.5
$ws.Range("L17").Value = 433831.5
$ws.Range("N17").Value = -434167.5
$ws.Range("H62").Value = 13499.889
$ws.Range("I62").Value = 12875
$ws.Range("J62").Value = 13999.8
$ws.Range("K62").Value = 12875
$ws.Range("L62").Value = 13999.8
$ws.Range("M62").Value = -12251
$ws.Range("N62").Value = -15247.8
$ws.Range("H65").Value = 13499.889
$ws.Range("I65").Value = 12875
$ws.Range("J65").Value = 13999.8
$ws.Range("K65").Value = 64375
$ws.Range("L65").Value = 69999
$ws.Range("M65").Value = -61255
$ws.Range("N65").Value = -76239
$ws.Range("H103").Value = 928.3333
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 928.3333
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2784.9999
$ws.Range("N103").Value = -3956.9999
$ws.Range("M103").ClearContents()
$ws.Range("H112").Value = 45629.15
$ws.Range("J112").Value = 29069.23
$ws.Range("L112").Value = 87207.69
$ws.Range("N112").Value = -89423.69
$ws.Range("H137").Value = 2206
$ws.Range("I137").Value = 1207.125
$ws.Range("K137").Value = 3621.375
$ws.Range("M137").Value = -1071.375
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 100010570
$ws.Range("I74").Value = 111122630
$ws.Range("K74").Value = 111122630
$ws.Range("M74").Value = -111121756
$ws.Range("H77").Value = 100010570
$ws.Range("I77").Value = 111122630
$ws.Range("K77").Value = 555613150
$ws.Range("M77").Value = -555608782
$ws.Range("H97").Value = 601.94446
$ws.Range("I97").Value = 599.4286
$ws.Range("J97").Value = 610.75
$ws.Range("K97").Value = 599.4286
$ws.Range("L97").Value = 610.75
$ws.Range("M97").Value = -103.4286
$ws.Range("N97").Value = -1602.75
$ws.Range("H122").Value = 1349.1389
$ws.Range("I122").Value = 1024.3871
$ws.Range("K122").Value = 3073.1613
$ws.Range("M122").Value = -623.1612999999998
$ws.Range("H132").Value = 2945167.5
$ws.Range("I132").Value = 4548782
$ws.Range("J132").Value = 5207.8335
$ws.Range("K132").Value = 13646346
$ws.Range("L132").Value = 15623.5005
$ws.Range("M132").Value = -13643816
$ws.Range("N132").Value = -20683.5005

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15069.667
$ws.Range("I26").Value = 15069.667
$ws.Range("K26").Value = 15069.667
$ws.Range("M26").Value = -14777.667
$ws.Range("H86").Value = 2025.9286
$ws.Range("I86").Value = 2008.5
$ws.Range("K86").Value = 2008.5
$ws.Range("M86").Value = -885.5
$ws.Range("H89").Value = 2025.9286
$ws.Range("I89").Value = 2008.5
$ws.Range("K89").Value = 10042.5
$ws.Range("M89").Value = -4426.5
$ws.Range("H100").Value = 21107.25
$ws.Range("I100").Value = 8000
$ws.Range("J100").Value = 25476.334
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 25476.334
$ws.Range("M100").Value = -6918
$ws.Range("N100").Value = -27640.334

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 13662.5
$ws.Range("J4").Value = 13662.5
$ws.Range("L4").Value = 13662.5
$ws.Range("N4").Value = -13886.5
$ws.Range("H7").Value = 232.41667
$ws.Range("I7").Value = 68.5
$ws.Range("J7").Value = 396.33334
$ws.Range("K7").Value = 68.5
$ws.Range("L7").Value = 396.33334
$ws.Range("M7").Value = 44.5
$ws.Range("N7").Value = -622.33334
$ws.Range("H28").Value = 5000
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5490
$ws.Range("H31").Value = 15043.857
$ws.Range("I31").Value = 10615.125
$ws.Range("K31").Value = 10615.125
$ws.Range("M31").Value = -10320.125
$ws.Range("H34").Value = 15043.857
$ws.Range("I34").Value = 10615.125
$ws.Range("K34").Value = 10615.125
$ws.Range("M34").Value = -10413.125
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 166670160
$ws.Range("I132").Value = 166670160
$ws.Range("K132").Value = 500010480
$ws.Range("M132").Value = -500007950

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3973887.8
$ws.Range("I4").Value = 4242619.5
$ws.Range("J4").Value = 2860570.8
$ws.Range("K4").Value = 12727858.5
$ws.Range("L4").Value = 8581712.399999999
$ws.Range("M4").Value = -12727746.5
$ws.Range("N4").Value = -8581936.399999999
$ws.Range("H5").Value = 111716.336
$ws.Range("J5").Value = 1183.3334
$ws.Range("L5").Value = 3550.0002
$ws.Range("N5").Value = -3774.0002
$ws.Range("H12").Value = 234.33333
$ws.Range("I12").Value = 151
$ws.Range("K12").Value = 453
$ws.Range("M12").Value = -280
$ws.Range("H17").Value = 641.1111
$ws.Range("I17").Value = 638.2857
$ws.Range("K17").Value = 1914.8571
$ws.Range("M17").Value = -1745.8571
$ws.Range("H32").Value = 2850
$ws.Range("J32").Value = 2850
$ws.Range("L32").Value = 8550
$ws.Range("N32").Value = -9116
$ws.Range("H113").Value = 112005.11
$ws.Range("I113").Value = 334050.34
$ws.Range("J113").Value = 982.5
$ws.Range("K113").Value = 1002151.02
$ws.Range("L113").Value = 2947.5
$ws.Range("M113").Value = -999981.02
$ws.Range("N113").Value = -7287.5
$ws.Range("H135").Value = 111716.336
$ws.Range("J135").Value = 1183.3334
$ws.Range("L135").Value = 10650.0006
$ws.Range("N135").Value = -15720.0006

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6885.3887
$ws.Range("I122").Value = 4448.909
$ws.Range("K122").Value = 13346.727
$ws.Range("M122").Value = -10896.727

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 29998.75
$ws.Range("J2").Value = 29998.75
$ws.Range("L2").Value = 29998.75
$ws.Range("N2").Value = -30222.75
$ws.Range("H93").Value = 3664.1667
$ws.Range("I93").Value = 3897
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 3897
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -2649
$ws.Range("N93").Value = -4996

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15784
$ws.Range("I2").Value = 10098.8
$ws.Range("J2").Value = 29997
$ws.Range("K2").Value = 10098.8
$ws.Range("L2").Value = 29997
$ws.Range("M2").Value = -9986.799999999999
$ws.Range("N2").Value = -30221
$ws.Range("H81").Value = 1499.2354
$ws.Range("I81").Value = 1565.8
$ws.Range("K81").Value = 3131.6
$ws.Range("M81").Value = -2070.6
$ws.Range("H84").Value = 1499.2354
$ws.Range("I84").Value = 1565.8
$ws.Range("K84").Value = 15658
$ws.Range("M84").Value = -10354
